$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capitalize the first letter of each PHENOTYPE (column A) value, rows 2-27
$phenotypes = @{
    2  = "Endometriosis"
    3  = "Endometriosis"
    4  = "Polycystic ovary syndrome"
    5  = "Polycystic ovary syndrome"
    6  = "Polycystic ovary syndrome"
    7  = "Polycystic ovary syndrome"
    8  = "Recurrent spontaneous abortion"
    9  = "Recurrent spontaneous abortion"
    10 = "Recurrent spontaneous abortion"
    11 = "Recurrent spontaneous abortion"
    12 = "Recurrent spontaneous abortion"
    13 = "Recurrent spontaneous abortion"
    14 = "DNA damage-related male infertility"
    15 = "DNA damage-related male infertility"
    16 = "Folic acid metabolism-related male infertility"
    17 = "Folic acid metabolism-related male infertility"
    18 = "Folic acid metabolism-related male infertility"
    19 = "Folic acid metabolism-related male infertility"
    20 = "Male infertility due to oxidative stress"
    21 = "Male infertility due to oxidative stress"
    22 = "Oligoasthenoteratozoospermia"
    23 = "Oligoasthenoteratozoospermia"
    24 = "Oligoasthenoteratozoospermia"
    25 = "Oligoasthenoteratozoospermia"
    26 = "Non-obstructive azoospermia"
    27 = "Non-obstructive azoospermia"
}

foreach ($row in $phenotypes.Keys) {
    $ws.Range("A$row").Value = $phenotypes[$row]
}

# Update the selected cell to match the saved view state (active cell A27)
$ws.Range("A27").Select()
